$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 3 (A3 = "H")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 185
$wsOff.Range("C3").Value = 132
$wsOff.Range("D3").Value = 58
$wsOff.Range("E3").Value = 27

# Sheet "DEF" - row 3 (A3 = "H")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 229
$wsDef.Range("C3").Value = 173
$wsDef.Range("D3").Value = 54
$wsDef.Range("E3").Value = 28
$wsDef.Range("F3").Value = 4
